$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "23×16=368"
$t.Cell(1,2).Range.Text = "14×64=896"
$t.Cell(1,3).Range.Text = "66×89=5874"
$t.Cell(1,4).Range.Text = "99×74=7326"
$t.Cell(1,5).Range.Text = "84×65=5460"
$t.Cell(5,1).Range.Text = "45×98=4410"
$t.Cell(5,2).Range.Text = "55×65=3575"
$t.Cell(5,3).Range.Text = "95×38=3610"
$t.Cell(5,4).Range.Text = "13×98=1274"
$t.Cell(5,5).Range.Text = "24×68=1632"
$t.Cell(10,1).Range.Text = "14×50=700"
$t.Cell(10,2).Range.Text = "93×87=8091"
$t.Cell(10,3).Range.Text = "63×73=4599"
$t.Cell(10,4).Range.Text = "62×44=2728"
$t.Cell(10,5).Range.Text = "62×89=5518"
$t.Cell(15,1).Range.Text = "80×73=5840"
$t.Cell(15,2).Range.Text = "34×37=1258"
$t.Cell(15,3).Range.Text = "82×44=3608"
$t.Cell(15,4).Range.Text = "95×89=8455"
$t.Cell(15,5).Range.Text = "25×54=1350"
$t.Cell(20,1).Range.Text = "68×93=6324"
$t.Cell(20,2).Range.Text = "94×23=2162"
$t.Cell(20,3).Range.Text = "98×50=4900"
$t.Cell(20,4).Range.Text = "80×44=3520"
$t.Cell(20,5).Range.Text = "92×74=6808"
